$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "use one stack to collect elements, for any pop/peek operation, move everything to second stack to get/remove top element, then move everything back.`nThen for any push operation, push to first stack`nEmpty is easy (call on first stack)"
$ws.Range("F2").Value = "push: make sure everything from s2 is in s1, then push to s1`npop: move everything to s2. now top of stack is first in queue`npeek: same as pop!`nEmpty: both are empty"
$ws.Range("G2").Value = "push: O(N)`npop: O(N)`npeek: O(N)`nEmpty: O(1)"
$ws.Range("H2").Value = "push: O(1)`npop: O(1)`npeek: O(1)`nEmpty: O(1)"

$ws.Range("E2:H2").Style = $ws.Range("C2").Style
$ws.Range("E2:H2").WrapText = $true

$ws.Rows.Item(2).RowHeight = 272

$excel.ActiveWindow.TopLeftCell = $ws.Range("B1")
$ws.Range("H3").Select()
